$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.290.68'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.864.21'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '234.09'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.4694'
$ws.Range("D7").Style = 'Normal'
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.2864'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.06555'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '21.70'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07868'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '96.15'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("D13").Value = '1.871.04'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '0.6933'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +1.85%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '5.090'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '265.95'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = '30.266.62'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '14.05'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.000007648'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +3.33%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = '2.115.53'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '5.229'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '6.187'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '9.381'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '167.26'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '18.80'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '1.938'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '0.09891'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.353'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '4.359'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.459'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '4.055'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.04745'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.131'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.7001'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.726'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.01869'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '2.792'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +5.85%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '6.195'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '72.39'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '1.940'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.8429'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.4163'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '102.38'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -0.89%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '7.106'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '941.89'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '9.137'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '34.50'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.05679'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.37%  '
